# Se agrega al plan de testeo el modulo de hoja de ruta
# Adds the "Hoja de Ruta" (routing sheet) test-plan module to sheet 1 ("Wildo"),
# mirroring the existing module sections already present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wildo")
$ws.Activate()

# --- Footer / signature row (row 194) ---------------------------------
$ws.Cells.Item(194, 2).Value = "Programador: Wildo Monges"
$ws.Cells.Item(194, 3).Value = "Tester: Leois linka"

# --- New module header (row 195) : bold, matches "Modulo" header style -
$ws.Cells.Item(195, 1).Value = "Modulo"
$ws.Cells.Item(195, 2).Value = "Hoja de Ruta(Tener en cuenta el orden de los productos que se van registrando. El primero"
$ws.Cells.Item(195, 3).Value = "debe estar en la primera posicion, el segundo en el segundo y asi….)"
$ws.Range("A195:C195").Font.Bold = $true

# --- Test case 1 (rows 196-198) ----------------------------------------
$ws.Cells.Item(196, 2).Value = "1) Guardar una hoja de ruta sin detalles"
$ws.Cells.Item(196, 3).Value = "El boton guardar no debe habilitarse, en caso de hackear el javascript"
$ws.Cells.Item(197, 3).Value = "del boton guardar el servidor muestra un mensaje:"
$ws.Cells.Item(198, 3).Value = '"Prohibido guardar sin agregar algun producto..", validado en el servidor'

# --- Test case 2 (rows 200-201) ----------------------------------------
$ws.Cells.Item(200, 2).Value = '2)En el campo codigo: ingresar 10 codigos de barra de productos en estado="No Enviado"'
$ws.Cells.Item(200, 3).Value = "Debe aparecer en la lista de los detalles los 10 productos en el orden"
$ws.Cells.Item(201, 3).Value = "que se fueron ingresando en que hayas ingresado."

# --- Test case 3 (rows 203-204) ----------------------------------------
$ws.Cells.Item(203, 2).Value = "3) Guardar la hoja de ruta despues de haber seleccionado la zona y agregado los detalles"
$ws.Cells.Item(203, 3).Value = "Debe redireccionarte a la pagina show.html donde se muestra la hoja de "
$ws.Cells.Item(204, 3).Value = "ruta creada"

# --- Test case 4 (row 205) ----------------------------------------------
$ws.Cells.Item(205, 2).Value = "4) Hacer click en Imprimir la hoja de ruta"
$ws.Cells.Item(205, 3).Value = "Debe generar como pdf para imprimir la hoja de ruta"

# --- Test case 5 (row 207) ----------------------------------------------
$ws.Cells.Item(207, 2).Value = "5) Click en Regresar"
$ws.Cells.Item(207, 3).Value = "Debe redireccionar a la pagina de buscar hoja de ruta"

# --- Test case 6 (row 209) ----------------------------------------------
$ws.Cells.Item(209, 2).Value = "6) Agregar un codigo de barra ya existente de un producto registrado"
$ws.Cells.Item(209, 3).Value = "No debe de aparecer en la lista"

# --- Test case 7 (row 210) ----------------------------------------------
$ws.Cells.Item(210, 2).Value = "7) Click en eliminar el producto de los detalles de hoja de ruta al estar creando la hoja"
$ws.Cells.Item(210, 3).Value = "Debe eliminar el producto de la lista"

# --- Test case 8 (row 211) ----------------------------------------------
$ws.Cells.Item(211, 2).Value = "8) Agregar el producto anteriormente eliminado"
$ws.Cells.Item(211, 3).Value = "Debe agregarse a la lista de los detalles"

# --- Test case 9 (row 212) ----------------------------------------------
$ws.Cells.Item(212, 2).Value = '9) Agregar el texto "Estos productos son de carácter privado" en el campo comentario'
$ws.Cells.Item(212, 3).Value = "Al hacer click sobre imprimir debe aparecer en el pdf a lado del campo OBS:"

# --- Test case 10 (rows 213-214) ----------------------------------------
$ws.Cells.Item(213, 2).Value = "10) Ir a Buscar Hoja de Ruta, buscar una hoja de ruta entre un rango de fechas y cuyo"
$ws.Cells.Item(214, 2).Value = "estado sea procesado, e intentar editar"
$ws.Cells.Item(214, 3).Value = "El boton edit debe estar desactivado si la hoja de ruta esta procesada"

# --- Test case 11 (rows 215-216) ----------------------------------------
$ws.Cells.Item(215, 2).Value = '11) Buscar una hoja de ruta con estado "En proceso", seleccionar e intentar editar'
$ws.Cells.Item(215, 3).Value = "Debe redireccionarte a la interfaz para editar la hoja e ruta con los campos"
$ws.Cells.Item(216, 3).Value = "de la hoja de ruta seteadas"

# --- Test case 12 (rows 217-218) ----------------------------------------
$ws.Cells.Item(217, 2).Value = "12) Al editar la hoja de ruta, eliminar alguno de sus productos y guardar"
$ws.Cells.Item(217, 3).Value = 'Esto productos eliminados deben cambiar su estado de "Enviado" a'
$ws.Cells.Item(218, 3).Value = '"No Enviado", y la hoja de ruta ahora debe aparecer sin estos productos'

# --- Update the view: scroll near the bottom and select C194 -----------
$ws.Range("C194").Select()
